$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three data rows (2-4) were appended again as rows 5-7 (an exact
# duplicate of the existing match-log entries), growing the used range
# from A1:K4 to A1:K7. Copy preserves the original text formatting
# (values stored as text, including the non-breaking space in
# "Navdeep Saini ") exactly as in the source rows.
$src = $ws.Range("A2:K4")
$dst = $ws.Range("A5:K7")
$src.Copy($dst)
